$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 32   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/10/2025  Through  11/16/2025"

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("L15").Value = -20
$ws.Range("N15").Value = -60
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -41.176470588235
$ws.Range("J16").Value = 106
$ws.Range("K16").Value = -16.981132075471
$ws.Range("L16").Value = 6.024096385542
$ws.Range("N16").Value = -83.489681050656
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 11
$ws.Range("H17").Value = 10
$ws.Range("I17").Value = 82
$ws.Range("J17").Value = 106
$ws.Range("K17").Value = -22.641509433962
$ws.Range("L17").Value = -14.583333333333
$ws.Range("M17").Value = 60.784313725490
$ws.Range("N17").Value = -3.529411764705
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -45.454545454545
$ws.Range("I18").Value = 118
$ws.Range("J18").Value = 123
$ws.Range("K18").Value = -4.065040650406
$ws.Range("L18").Value = 22.916666666666
$ws.Range("M18").Value = 20.408163265306
$ws.Range("N18").Value = -85.359801488833
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 23.076923076923
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 67
$ws.Range("H19").Value = -34.328358208955
$ws.Range("I19").Value = 647
$ws.Range("J19").Value = 694
$ws.Range("K19").Value = -6.772334293948
$ws.Range("L19").Value = -0.614439324116
$ws.Range("M19").Value = 11.551724137931
$ws.Range("N19").Value = -59.738643434972
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -75
$ws.Range("J20").Value = 52
$ws.Range("K20").Value = -23.076923076923
$ws.Range("L20").Value = -55.056179775280
$ws.Range("N20").Value = -96.078431372549
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -5
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -34.210526315789
$ws.Range("I21").Value = 983
$ws.Range("J21").Value = 1089
$ws.Range("K21").Value = -9.733700642791
$ws.Range("L21").Value = -4.470359572400
$ws.Range("M21").Value = 16.607354685646
$ws.Range("N21").Value = -75.859528487229
$ws.Range("L22").Value = -29.166666666666
$ws.Range("D23").Value = 2
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -62.5
$ws.Range("J23").Value = 41
$ws.Range("K23").Value = -29.268292682926
$ws.Range("M23").Value = 20.833333333333
$ws.Range("C24").Value = 32
$ws.Range("E24").Value = -5.882352941176
$ws.Range("G24").Value = 135
$ws.Range("H24").Value = -10.370370370370
$ws.Range("I24").Value = 1261
$ws.Range("J24").Value = 1194
$ws.Range("K24").Value = 5.611390284757
$ws.Range("L24").Value = 21.953578336557
$ws.Range("M24").Value = 34.434968017057
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = -3.703703703703
$ws.Range("F25").Value = 97
$ws.Range("G25").Value = 119
$ws.Range("H25").Value = -18.487394957983
$ws.Range("I25").Value = 929
$ws.Range("J25").Value = 964
$ws.Range("K25").Value = -3.630705394190
$ws.Range("L25").Value = 15.980024968789
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = -80
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = -10
$ws.Range("I26").Value = 218
$ws.Range("J26").Value = 205
$ws.Range("K26").Value = 6.341463414634
$ws.Range("L26").Value = 9.547738693467
$ws.Range("M26").Value = -12.096774193548
$ws.Range("L27").Value = -45
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 37
$ws.Range("J28").Value = 33
$ws.Range("K28").Value = 12.121212121212
$ws.Range("L28").Value = 5.714285714285
$ws.Range("J31").Value = 13
$ws.Range("K31").Value = -46.153846153846
$ws.Range("L31").Value = -56.25

# --- Value/type + style transitions (numeric <-> "no data" placeholder text) ---
# Stable, unedited reference cells used purely as formatting sources:
#   C14 = style for "0" placeholder (text, s=13, shared string "0")
#   E14 = style for "***.*" placeholder (text, s=13, shared string "***.*")
#   F15 = style for plain integer count (s=15)
#   L14 = style for percent/decimal number (s=14)

$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("C16").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("C18").Value = 2
$ws.Range("F15").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("C22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$ws.Range("D30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$ws.Range("D31").Value = 1
$ws.Range("F15").Copy()
$ws.Range("D31").PasteSpecial(-4122)

$ws.Range("E31").Value = -100
$ws.Range("L14").Copy()
$ws.Range("E31").PasteSpecial(-4122)

$ws.Range("G31").Value = 1
$ws.Range("F15").Copy()
$ws.Range("G31").PasteSpecial(-4122)

$ws.Range("H31").Value = -100
$ws.Range("L14").Copy()
$ws.Range("H31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

Write-Output "edit complete"
